$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1 (05:05 -> 05:35)
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 05:35"

# Update India row (row 14): Casos totales, Nuevos casos, Casos activos, Recuperados
$ws.Range("B14").Value = 131868
$ws.Range("C14").Value = 445
$ws.Range("D14").Value = 54441
$ws.Range("E14").Value = 73559

# Honduras moved up (new data pushed it above Hungria and Sudan).
# Row 72 used to be Hungria, row 73 Sudan, row 74 Honduras.
# New order: row 72 Honduras (updated data), row 73 Hungria (old row72 data), row 74 Sudan (old row73 data).
$ws.Range("A72").Value = "Honduras"
$ws.Range("B72").Value = 3743
$ws.Range("C72").Value = 266
$ws.Range("D72").Value = 455
$ws.Range("E72").Value = 3114
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 7
$ws.Range("H72").Value = 174

$ws.Range("A73").Value = "Hungria"
$ws.Range("B73").Value = 3713
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 1655
$ws.Range("E73").Value = 1576
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 482

$ws.Range("A74").Value = "Sudan"
$ws.Range("B74").Value = 3628
$ws.Range("C74").Value = 250
$ws.Range("D74").Value = 424
$ws.Range("E74").Value = 3058
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 9
$ws.Range("H74").Value = 146

# Update El Salvador row (row 92): Casos activos / Recuperados
$ws.Range("D92").Value = 574
$ws.Range("E92").Value = 1212

# Update Mongolia row (row 161): Casos activos / Recuperados
$ws.Range("D161").Value = 32
$ws.Range("E161").Value = 109
